$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# These cells already store date/time as plain text (shared-string) values,
# not numeric Excel dates, so a direct .Value assignment (without touching
# NumberFormat) keeps them text-typed and preserves their existing style.
$overview.Range("G2").Value = "2016-08-27 07:01:19"

$zhcn.Range("H2").Value = "2016-08-27 07:01:14"
$zhcn.Range("K2").Value = "2016-08-27 07:01:43"

$dede.Range("H2").Value = "2016-08-27 07:01:19"
$dede.Range("K2").Value = "2016-08-27 07:01:50"
